$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 = "Experimental" property; previously had no value in column B.
$ws.Range("B7").Value = "true"

# Row 8 = "Date" property; update the timestamp value.
$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"
